$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column (Price) cells: force text format so numeric-looking strings
# (e.g. "570.60", "3.50") are not coerced to numbers and lose trailing zeros,
# then restore the default "Normal" style so no stray formatting is left behind.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.916.38"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.42%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.406.18"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.19%  "
$ws.Range("E4").Value = "  -0.34%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.88%  "
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.527"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.386.57"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.90%  "
$ws.Range("E10").Value = "  -0.13%  "
$ws.Range("E11").Value = "  +0.45%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.07"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.338"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.42%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.96"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.58%  "
$ws.Range("E15").Value = "  -2.11%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.818.66"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.801.64"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.383.22"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.76%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "322.12"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.02"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.70%  "
$ws.Range("E23").Value = "  +1.31%  "
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("E25").Value = "  -7.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "64.40"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.53%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -8.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "578.61"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -8.03%  "
$ws.Range("E29").Value = "  -2.73%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0909"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.80"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.95%  "
$ws.Range("E33").Value = "  -2.87%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.130"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -8.50%  "
$ws.Range("E35").Value = "  +0.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.62"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.94%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.367"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.38"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.37%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "147.25"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.79%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.17"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.09"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.18%  "
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.65"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.67"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.96%  "
$ws.Range("E45").Value = "  -5.78%  "
$ws.Range("E46").Value = "  +19.64%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "140.34"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.80%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.50"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.35%  "
$ws.Range("E49").Value = "  -3.10%  "
$ws.Range("E50").Value = "  -4.73%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.36"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.01%  "
